$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.81136131286621
$ws.Range("D2").Value = -0.02863868713379247
$ws.Range("E2").Value = 0.0008201744007472505
$ws.Range("C3").Value = 29.75001335144043
$ws.Range("D3").Value = -0.05998664855957259
$ws.Range("E3").Value = 0.003598398005409672
$ws.Range("C4").Value = 29.76519584655762
$ws.Range("D4").Value = -0.1548041534423845
$ws.Range("E4").Value = 0.02396432592301333
$ws.Range("C5").Value = 29.79987907409668
$ws.Range("D5").Value = -0.1801209259033243
$ws.Range("E5").Value = 0.03244354794827084
$ws.Range("B6").Value = 30.03999999999999
$ws.Range("C6").Value = 30.01015281677246
$ws.Range("D6").Value = -0.0298471832275311
$ws.Range("E6").Value = 0.000890854346617814
$ws.Range("B7").Value = 30.21000000000001
$ws.Range("C7").Value = 30.06319236755371
$ws.Range("D7").Value = -0.146807632446297
$ws.Range("E7").Value = 0.02155248094448704
$ws.Range("C8").Value = 30.18989372253418
$ws.Range("D8").Value = -0.03010627746581918
$ws.Range("E8").Value = 0.0009063879428488914
$ws.Range("C9").Value = 30.25397300720215
$ws.Range("D9").Value = -0.126026992797847
$ws.Range("E9").Value = 0.01588280291366858
$ws.Range("C10").Value = 30.50796318054199
$ws.Range("D10").Value = 0.06796318054199446
$ws.Range("E10").Value = 0.004618993909383734
$ws.Range("C11").Value = 30.37946128845215
$ws.Range("D11").Value = -0.1005387115478555
$ws.Range("E11").Value = 0.0101080325197029
$ws.Range("C12").Value = 30.38317680358887
$ws.Range("D12").Value = -0.3068231964111305
$ws.Range("E12").Value = 0.09414047385594319
$ws.Range("C13").Value = 30.47823143005371
$ws.Range("D13").Value = -0.2717685699462891
$ws.Range("E13").Value = 0.07385815561065101
$ws.Range("C14").Value = 30.61132621765137
$ws.Range("D14").Value = -0.3286737823486305
$ws.Range("E14").Value = 0.108026455203355
$ws.Range("C15").Value = 30.73988914489746
$ws.Range("D15").Value = -0.2101108551025419
$ws.Range("E15").Value = 0.04414657143192136
$ws.Range("C16").Value = 31.15546226501465
$ws.Range("D16").Value = 0.1354622650146524
$ws.Range("E16").Value = 0.01835002524289992
$ws.Range("C17").Value = 31.3719310760498
$ws.Range("D17").Value = 0.2519310760498001
$ws.Range("E17").Value = 0.06346926707961018
$ws.Range("C18").Value = 31.45668601989746
$ws.Range("D18").Value = 0.1766860198974598
$ws.Range("E18").Value = 0.03121794962720556
$ws.Range("C19").Value = 31.35834312438965
$ws.Range("D19").Value = -0.02165687561034702
$ws.Range("E19").Value = 0.0004690202612020434
$ws.Range("C20").Value = 31.55037498474121
$ws.Range("D20").Value = -0.02962501525878736
$ws.Range("E20").Value = 0.0008776415290833837
$ws.Range("B21").Value = 31.65000000000001
$ws.Range("C21").Value = 32.06996536254883
$ws.Range("D21").Value = 0.4199653625488224
$ws.Range("E21").Value = 0.1763709057407639
$ws.Range("C22").Value = 32.6776237487793
$ws.Range("D22").Value = 0.7976237487793014
$ws.Range("E22").Value = 0.6362036446167462
$ws.Range("C23").Value = 32.63322830200195
$ws.Range("D23").Value = 0.353228302001952
$ws.Range("E23").Value = 0.1247702333351822
$ws.Range("C24").Value = 32.70923233032227
$ws.Range("D24").Value = 0.2592323303222628
$ws.Range("E24").Value = 0.06720140108431076
$ws.Range("B25").Value = 32.84999999999999
$ws.Range("C25").Value = 32.97509384155273
$ws.Range("D25").Value = 0.1250938415527401
$ws.Range("E25").Value = 0.01564846919442204
$ws.Range("B26").Value = 32.90000000000001
$ws.Range("C26").Value = 33.11442184448242
$ws.Range("D26").Value = 0.2144218444824162
$ws.Range("E26").Value = 0.04597672739124147
$ws.Range("B27").Value = 33.09999999999999
$ws.Range("C27").Value = 33.00744247436523
$ws.Range("D27").Value = -0.09255752563475994
$ws.Range("E27").Value = 0.008566895551629244
$ws.Range("B28").Value = 33.40000000000001
$ws.Range("C28").Value = 33.61306381225586
$ws.Range("D28").Value = 0.2130638122558537
$ws.Range("E28").Value = 0.04539618809299767
$ws.Range("C29").Value = 33.67734146118164
$ws.Range("D29").Value = -0.02265853881836222
$ws.Range("E29").Value = 0.0005134093813832275
$ws.Range("B30").Value = 34.09999999999999
$ws.Range("C30").Value = 33.8823127746582
$ws.Range("D30").Value = -0.2176872253417912
$ws.Range("E30").Value = 0.04738772807700778
$ws.Range("B31").Value = 34.40000000000001
$ws.Range("C31").Value = 34.47198104858398
$ws.Range("D31").Value = 0.07198104858397869
$ws.Range("E31").Value = 0.0051812713552491
$ws.Range("B32").Value = 34.90000000000001
$ws.Range("C32").Value = 35.08790969848633
$ws.Range("D32").Value = 0.1879096984863224
$ws.Range("E32").Value = 0.03531005478522061
$ws.Range("C33").Value = 35.6751823425293
$ws.Range("D33").Value = 0.3751823425292997
$ws.Range("E33").Value = 0.1407617901457728
$ws.Range("C34").Value = 35.94461441040039
$ws.Range("D34").Value = 0.2446144104003878
$ws.Range("E34").Value = 0.05983620977552934
$ws.Range("C35").Value = 36.12253189086914
$ws.Range("D35").Value = -0.1774681091308565
$ws.Range("E35").Value = 0.0314949297584816
$ws.Range("C36").Value = 36.61664962768555
$ws.Range("D36").Value = -0.1833503723144503
$ws.Range("E36").Value = 0.03361735902784754
$ws.Range("C37").Value = 37.1579704284668
$ws.Range("D37").Value = -0.1420295715332003
$ws.Range("E37").Value = 0.02017239918990446
$ws.Range("B38").Value = 37.90000000000001
$ws.Range("C38").Value = 37.8542594909668
$ws.Range("D38").Value = -0.04574050903320881
$ws.Range("E38").Value = 0.002092194166617057
$ws.Range("C39").Value = 38.33866500854492
$ws.Range("D39").Value = -0.1613349914550781
$ws.Range("E39").Value = 0.02602897946781013
$ws.Range("B40").Value = 38.90000000000001
$ws.Range("C40").Value = 39.00360488891602
$ws.Range("D40").Value = 0.1036048889160099
$ws.Range("E40").Value = 0.01073397300729876
$ws.Range("B41").Value = 39.40000000000001
$ws.Range("C41").Value = 39.52373123168945
$ws.Range("D41").Value = 0.1237312316894474
$ws.Range("E41").Value = 0.01530941769538772
$ws.Range("B42").Value = 39.90000000000001
$ws.Range("C42").Value = 39.62896347045898
$ws.Range("D42").Value = -0.2710365295410213
$ws.Range("E42").Value = 0.07346080034564091
$ws.Range("B43").Value = 40.09999999999999
$ws.Range("C43").Value = 39.96549606323242
$ws.Range("D43").Value = -0.1345039367675724
$ws.Range("E43").Value = 0.01809130900597513
$ws.Range("B44").Value = 40.59999999999999
$ws.Range("C44").Value = 40.37540817260742
$ws.Range("D44").Value = -0.2245918273925724
$ws.Range("E44").Value = 0.05044148893153505
$ws.Range("B45").Value = 40.90000000000001
$ws.Range("C45").Value = 40.5744743347168
$ws.Range("D45").Value = -0.3255256652832088
$ws.Range("E45").Value = 0.1059669587580757
$ws.Range("B46").Value = 41.20000000000001
$ws.Range("C46").Value = 41.21307373046875
$ws.Range("D46").Value = 0.01307373046874005
$ws.Range("E46").Value = 0.000170922428369262
$ws.Range("C47").Value = 41.28318023681641
$ws.Range("D47").Value = -0.2168197631835938
$ws.Range("E47").Value = 0.04701080970698968
$ws.Range("C48").Value = 41.82052230834961
$ws.Range("D48").Value = 0.02052230834961222
$ws.Range("E48").Value = 0.0004211651399965633
$ws.Range("C49").Value = 41.57979202270508
$ws.Range("D49").Value = -0.6202079772949247
$ws.Range("E49").Value = 0.3846579351002619
$ws.Range("C50").Value = 43.23081970214844
$ws.Range("D50").Value = 0.5308197021484347
$ws.Range("E50").Value = 0.2817695561889529
$ws.Range("C51").Value = 43.8011360168457
$ws.Range("D51").Value = 0.1011360168457074
$ws.Range("E51").Value = 0.01022849390341521
$ws.Range("C52").Value = -0.07380088806155527
$ws.Range("E52").Value = 3.070135179046037
$ws.Range("E53").Value = 0.06140270358092074

Write-Host "Applied partial model results"